$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processes")

$ws.Range("B2").Value = 3903488
$ws.Range("B3").Value = 4063232
$ws.Range("B4").Value = 3932160
$ws.Range("B5").Value = 62607360
$ws.Range("B6").Value = 16486400
$ws.Range("B7").Value = 21458944
$ws.Range("B8").Value = 50139136
$ws.Range("B9").Value = 110174208
$ws.Range("B10").Value = 14221312
